$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3871.0908
$ws.Range("I64").Value = 3884.6086
$ws.Range("J64").Value = 3840
$ws.Range("K64").Value = 3884.6086
$ws.Range("L64").Value = 3840
$ws.Range("M64").Value = -3636.6086
$ws.Range("N64").Value = -4336
$ws.Range("H67").Value = 3871.0908
$ws.Range("I67").Value = 3884.6086
$ws.Range("J67").Value = 3840
$ws.Range("K67").Value = 3884.6086
$ws.Range("L67").Value = 3840
$ws.Range("M67").Value = -3026.6086
$ws.Range("N67").Value = -5556
$ws.Range("H74").Value = 4327.7856
$ws.Range("I74").Value = 4170
$ws.Range("J74").Value = 4722.25
$ws.Range("K74").Value = 4170
$ws.Range("L74").Value = 4722.25
$ws.Range("M74").Value = -3234
$ws.Range("N74").Value = -6594.25
$ws.Range("H76").Value = 8338.619000000001
$ws.Range("I76").Value = 11241.917
$ws.Range("J76").Value = 4467.5557
$ws.Range("K76").Value = 11241.917
$ws.Range("L76").Value = 4467.5557
$ws.Range("M76").Value = -10926.917
$ws.Range("N76").Value = -5097.5557
$ws.Range("H77").Value = 4327.7856
$ws.Range("I77").Value = 4170
$ws.Range("J77").Value = 4722.25
$ws.Range("K77").Value = 20850
$ws.Range("L77").Value = 23611.25
$ws.Range("M77").Value = -16170
$ws.Range("N77").Value = -32971.25
$ws.Range("H79").Value = 8338.619000000001
$ws.Range("I79").Value = 11241.917
$ws.Range("J79").Value = 4467.5557
$ws.Range("K79").Value = 11241.917
$ws.Range("L79").Value = 4467.5557
$ws.Range("M79").Value = -10149.917
$ws.Range("N79").Value = -6651.5557
$ws.Range("H93").Value = 22146.285
$ws.Range("J93").Value = 22146.285
$ws.Range("L93").Value = 22146.285
$ws.Range("N93").Value = -27138.285
$ws.Range("H98").Value = 3962.875
$ws.Range("I98").Value = 2570.077
$ws.Range("J98").Value = 9998.333000000001
$ws.Range("K98").Value = 2570.077
$ws.Range("L98").Value = 9998.333000000001
$ws.Range("M98").Value = -1072.077
$ws.Range("N98").Value = -12994.333
$ws.Range("H122").Value = 3962.875
$ws.Range("I122").Value = 2570.077
$ws.Range("J122").Value = 9998.333000000001
$ws.Range("K122").Value = 7710.231000000001
$ws.Range("L122").Value = 29994.999
$ws.Range("M122").Value = -5260.231000000001
$ws.Range("N122").Value = -34894.999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5573.7534
$ws.Range("I32").Value = 4548.8213
$ws.Range("J32").Value = 9988.846
$ws.Range("K32").Value = 4548.8213
$ws.Range("L32").Value = 9988.846
$ws.Range("M32").Value = -4261.8213
$ws.Range("N32").Value = -10562.846
$ws.Range("H45").Value = 4874.7036
$ws.Range("I45").Value = 7339.8125
$ws.Range("J45").Value = 1289.091
$ws.Range("K45").Value = 7339.8125
$ws.Range("L45").Value = 1289.091
$ws.Range("M45").Value = -6962.8125
$ws.Range("N45").Value = -2043.091
$ws.Range("H61").Value = 6854.476
$ws.Range("I61").Value = 7370.737
$ws.Range("K61").Value = 7370.737
$ws.Range("M61").Value = -7158.737
$ws.Range("H74").Value = 1746
$ws.Range("I74").Value = 1660.5186
$ws.Range("K74").Value = 1660.5186
$ws.Range("M74").Value = -786.5186000000001
$ws.Range("H77").Value = 1746
$ws.Range("I77").Value = 1660.5186
$ws.Range("K77").Value = 8302.593000000001
$ws.Range("M77").Value = -3934.593000000001
$ws.Range("H132").Value = 3442.641
$ws.Range("I132").Value = 1794.3334
$ws.Range("J132").Value = 4855.476
$ws.Range("K132").Value = 5383.0002
$ws.Range("L132").Value = 14566.428
$ws.Range("M132").Value = -2853.0002
$ws.Range("N132").Value = -19626.428
$ws.Range("H135").Value = 85000
$ws.Range("J135").Value = 85000
$ws.Range("L135").Value = 85000
$ws.Range("N135").Value = -95140
$ws.Range("H136").Value = 6854.476
$ws.Range("I136").Value = 7370.737
$ws.Range("K136").Value = 22112.211
$ws.Range("M136").Value = -19562.211
$ws.Range("H141").Value = 37100
$ws.Range("J141").Value = 37100
$ws.Range("L141").Value = 37100
$ws.Range("N141").Value = -47460

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4578.6587
$ws.Range("I134").Value = 6197.1665
$ws.Range("J134").Value = 2293.7058
$ws.Range("K134").Value = 18591.4995
$ws.Range("L134").Value = 6881.117400000001
$ws.Range("M134").Value = -16056.4995
$ws.Range("N134").Value = -11951.1174

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8548140
$ws.Range("I16").Value = 15385674
$ws.Range("J16").Value = 1222.5
$ws.Range("K16").Value = 15385674
$ws.Range("L16").Value = 1222.5
$ws.Range("M16").Value = -15385387
$ws.Range("N16").Value = -1796.5
$ws.Range("H62").Value = 6022.778
$ws.Range("I62").Value = 6315
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 6315
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -5691
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 6022.778
$ws.Range("I65").Value = 6315
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 31575
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -28455
$ws.Range("N65").Value = -31240
$ws.Range("H105").Value = 37038984
$ws.Range("I105").Value = 37038984
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 37038984
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -37037237
$ws.Range("H113").Value = 8548140
$ws.Range("I113").Value = 15385674
$ws.Range("J113").Value = 1222.5
$ws.Range("K113").Value = 15385674
$ws.Range("L113").Value = 1222.5
$ws.Range("M113").Value = -15383504
$ws.Range("N113").Value = -5562.5
$ws.Range("N105").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1801.75
$ws.Range("I68").Value = 442.8
$ws.Range("J68").Value = 4066.6667
$ws.Range("K68").Value = 1328.4
$ws.Range("L68").Value = 12200.0001
$ws.Range("M68").Value = -517.4000000000001
$ws.Range("N68").Value = -13822.0001
$ws.Range("H71").Value = 1801.75
$ws.Range("I71").Value = 442.8
$ws.Range("J71").Value = 4066.6667
$ws.Range("K71").Value = 3985.2
$ws.Range("L71").Value = 36600.0003
$ws.Range("M71").Value = 70.79999999999973
$ws.Range("N71").Value = -44712.0003
$ws.Range("H132").Value = 1795515.5
$ws.Range("J132").Value = 1824933.8
$ws.Range("L132").Value = 16424404.2
$ws.Range("N132").Value = -16429464.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5952.5366
$ws.Range("I70").Value = 5972.484
$ws.Range("J70").Value = 5890.7
$ws.Range("K70").Value = 5972.484
$ws.Range("L70").Value = 5890.7
$ws.Range("M70").Value = -5702.484
$ws.Range("N70").Value = -6430.7
$ws.Range("H73").Value = 5952.5366
$ws.Range("I73").Value = 5972.484
$ws.Range("J73").Value = 5890.7
$ws.Range("K73").Value = 5972.484
$ws.Range("L73").Value = 5890.7
$ws.Range("M73").Value = -5036.484
$ws.Range("N73").Value = -7762.7
$ws.Range("H80").Value = 2499.9844
$ws.Range("I80").Value = 2499.9844
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2499.9844
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1501.9844
$ws.Range("H83").Value = 2499.9844
$ws.Range("I83").Value = 2499.9844
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 12499.922
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -7507.921999999999
$ws.Range("H126").Value = 5866.24
$ws.Range("I126").Value = 7013.421
$ws.Range("J126").Value = 2233.5
$ws.Range("K126").Value = 21040.263
$ws.Range("L126").Value = 6700.5
$ws.Range("M126").Value = -18570.263
$ws.Range("N126").Value = -11640.5
$ws.Range("H132").Value = 4157.92
$ws.Range("I132").Value = 6028.8
$ws.Range("J132").Value = 2910.6667
$ws.Range("K132").Value = 18086.4
$ws.Range("L132").Value = 8732.000100000001
$ws.Range("M132").Value = -15556.4
$ws.Range("N132").Value = -13792.0001
$ws.Range("N80").ClearContents()
$ws.Range("N83").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 253750
$ws.Range("I18").Value = 500500
$ws.Range("K18").Value = 500500
$ws.Range("M18").Value = -500328
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("H132").Value = 26970508
$ws.Range("I132").Value = 32748766
$ws.Range("J132").Value = 5301.3335
$ws.Range("K132").Value = 98246298
$ws.Range("L132").Value = 15904.0005
$ws.Range("M132").Value = -98243768
$ws.Range("N132").Value = -20964.0005
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H136").Value = 4448.683
$ws.Range("I136").Value = 4935.4517
$ws.Range("K136").Value = 14806.3551
$ws.Range("M136").Value = -12256.3551
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H140").Value = 46273
$ws.Range("J140").Value = 46273
$ws.Range("L140").Value = 46273
$ws.Range("N140").Value = -56633
$ws.Range("H141").Value = 48057.145
$ws.Range("J141").Value = 48057.145
$ws.Range("L141").Value = 48057.145
$ws.Range("N141").Value = -58417.145
$ws.Range("N20").ClearContents()
$ws.Range("N135").ClearContents()
$ws.Range("N138").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1891.1538
$ws.Range("I113").Value = 1177.875
$ws.Range("J113").Value = 3032.4
$ws.Range("K113").Value = 3533.625
$ws.Range("L113").Value = 9097.200000000001
$ws.Range("M113").Value = -1363.625
$ws.Range("N113").Value = -13437.2
$ws.Range("H136").Value = 5358.5713
$ws.Range("I136").Value = 17266.666
$ws.Range("J136").Value = 2110.9092
$ws.Range("K136").Value = 51799.99800000001
$ws.Range("L136").Value = 6332.7276
$ws.Range("M136").Value = -49249.99800000001
$ws.Range("N136").Value = -11432.7276
$ws.Range("H139").Value = 59800
$ws.Range("J139").Value = 59800
$ws.Range("L139").Value = 59800
$ws.Range("N139").Value = -70080
